# Apply cryptocurrency price/volume updates described in the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.767.82'
$ws.Range('E2').Value = '  +4.76%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.982.25'
$ws.Range('E3').Value = '  +2.74%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.20%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.55'
$ws.Range('E5').Value = '  +2.08%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.42'
$ws.Range('E6').Value = '  +6.39%  '

$ws.Range('E7').Value = '  +0.06%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.977.19'
$ws.Range('E8').Value = '  +2.59%  '

$ws.Range('E9').Value = '  +0.12%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.95'
$ws.Range('E10').Value = '  +5.10%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.153'
$ws.Range('E11').Value = '  +3.29%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.446'
$ws.Range('E12').Value = '  +2.32%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000239'
$ws.Range('E13').Value = '  +3.00%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.18'
$ws.Range('E14').Value = '  +6.23%  '

$ws.Range('E15').Value = '  +0.72%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.857.63'
$ws.Range('E16').Value = '  +4.96%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.483.56'
$ws.Range('E17').Value = '  +2.88%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.89'
$ws.Range('E18').Value = '  +3.41%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.984.24'
$ws.Range('E19').Value = '  +3.81%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '446.70'
$ws.Range('E20').Value = '  +2.27%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.62'
$ws.Range('E21').Value = '  +2.27%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.677'
$ws.Range('E22').Value = '  +2.89%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.29'
$ws.Range('E23').Value = '  +5.26%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '81.03'
$ws.Range('E24').Value = '  +1.80%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.88'
$ws.Range('E25').Value = '  +7.02%  '

$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.23'
$ws.Range('E26').Value = '  +3.08%  '

$ws.Range('B27').Value = 'Fetch.AI'
$ws.Range('C27').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.18'
$ws.Range('E27').Value = '  +6.98%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.37'
$ws.Range('E29').Value = '  +14.13%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.70'
$ws.Range('E30').Value = '  +8.43%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0000107'
$ws.Range('E31').Value = '  +3.16%  '

$ws.Range('E32').Value = '  +2.28%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.110'
$ws.Range('E33').Value = '  +2.57%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.55'
$ws.Range('E34').Value = '  +3.39%  '

$ws.Range('E35').Value = '  +0.19%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.982'
$ws.Range('E36').Value = '  +2.36%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.64'
$ws.Range('E37').Value = '  +3.42%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.11'
$ws.Range('E38').Value = '  +7.75%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.98'
$ws.Range('E39').Value = '  +5.79%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '48.88'
$ws.Range('E40').Value = '  -0.22%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '43.88'
$ws.Range('E41').Value = '  +10.46%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.120'
$ws.Range('E42').Value = '  +3.98%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.294'
$ws.Range('E43').Value = '  +8.99%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.37'
$ws.Range('E44').Value = '  +0.94%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '383.70'
$ws.Range('E45').Value = '  +13.29%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.785.91'
$ws.Range('E46').Value = '  +3.14%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0347'
$ws.Range('E47').Value = '  +3.80%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '134.98'
$ws.Range('E48').Value = '  +0.74%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000220'
$ws.Range('E50').Value = '  +13.79%  '

$ws.Range('E51').Value = '  +1.50%  '

